$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54; this shifts rows 54..126 down to 55..127
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new record's data.
# Columns A,B,C,E,F,G,H,I,N,Q,R are constant across every data row in this sheet.
$ws.Cells.Item(54, 1).Value = 5
$ws.Cells.Item(54, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(54, 3).Value = "Maule"
$ws.Cells.Item(54, 4).Value = 45210
$ws.Cells.Item(54, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(54, 5).Value = 7
$ws.Cells.Item(54, 6).Value = 100112026
$ws.Cells.Item(54, 7).Value = "Haba"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 200
$ws.Cells.Item(54, 11).Value = 10000
$ws.Cells.Item(54, 12).Value = 10000
$ws.Cells.Item(54, 13).Value = 10000
$ws.Cells.Item(54, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(54, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(54, 16).Value = 400
$ws.Cells.Item(54, 17).Value = 25
$ws.Cells.Item(54, 18).Value = "Hortaliza"
